# Auto-generated edit script: refresh market-price-derived Leve profit data
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3945.3845
$ws.Range("J64").Value = 3899.5715
$ws.Range("L64").Value = 3899.5715
$ws.Range("N64").Value = -4395.5715
$ws.Range("H67").Value = 3945.3845
$ws.Range("J67").Value = 3899.5715
$ws.Range("L67").Value = 3899.5715
$ws.Range("N67").Value = -5615.5715
$ws.Range("H74").Value = 6969.7393
$ws.Range("I74").Value = 6209.7144
$ws.Range("K74").Value = 6209.7144
$ws.Range("M74").Value = -5273.7144
$ws.Range("H77").Value = 6969.7393
$ws.Range("I77").Value = 6209.7144
$ws.Range("K77").Value = 31048.572
$ws.Range("M77").Value = -26368.572
$ws.Range("H86").Value = 6999
$ws.Range("J86").Value = 6999
$ws.Range("L86").Value = 6999
$ws.Range("N86").Value = -9245
$ws.Range("H89").Value = 6999
$ws.Range("J89").Value = 6999
$ws.Range("L89").Value = 34995
$ws.Range("N89").Value = -46227
$ws.Range("H96").Value = 808.75
$ws.Range("I96").Value = 819.25
$ws.Range("J96").Value = 798.25
$ws.Range("K96").Value = 2457.75
$ws.Range("L96").Value = 2394.75
$ws.Range("M96").Value = -1084.75
$ws.Range("N96").Value = -5140.75
$ws.Range("H97").Value = 7021
$ws.Range("J97").Value = 7021
$ws.Range("L97").Value = 21063
$ws.Range("N97").Value = -22055
$ws.Range("H101").Value = 3443.1428
$ws.Range("I101").Value = 1534.2222
$ws.Range("J101").Value = 6879.2
$ws.Range("K101").Value = 4602.6666
$ws.Range("L101").Value = 20637.6
$ws.Range("M101").Value = -2980.6666
$ws.Range("N101").Value = -23881.6
$ws.Range("H134").Value = 42996
$ws.Range("J134").Value = 42996
$ws.Range("L134").Value = 42996
$ws.Range("N134").Value = -53136
$ws.Range("H135").Value = 1383.0555
$ws.Range("J135").Value = 1999.5
$ws.Range("L135").Value = 17995.5
$ws.Range("N135").Value = -23065.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1899
$ws.Range("I63").Value = 1899
$ws.Range("K63").Value = 1899
$ws.Range("M63").Value = -1213
$ws.Range("H66").Value = 1899
$ws.Range("I66").Value = 1899
$ws.Range("K66").Value = 9495
$ws.Range("M66").Value = -6063
$ws.Range("H88").Value = 26500
$ws.Range("H91").Value = 26500
$ws.Range("H132").Value = 3191.4
$ws.Range("I132").Value = 2379.2727
$ws.Range("K132").Value = 7137.8181
$ws.Range("M132").Value = -4607.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 22500
$ws.Range("J76").Value = 22500
$ws.Range("L76").Value = 22500
$ws.Range("N76").Value = -23130
$ws.Range("H79").Value = 22500
$ws.Range("J79").Value = 22500
$ws.Range("L79").Value = 22500
$ws.Range("N79").Value = -24684
$ws.Range("H86").Value = 6830.4346
$ws.Range("I86").Value = 2025.1765
$ws.Range("K86").Value = 2025.1765
$ws.Range("M86").Value = -902.1765
$ws.Range("H89").Value = 6830.4346
$ws.Range("I89").Value = 2025.1765
$ws.Range("K89").Value = 10125.8825
$ws.Range("M89").Value = -4509.8825
$ws.Range("H133").Value = 85000
$ws.Range("J133").Value = 85000
$ws.Range("L133").Value = 85000
$ws.Range("N133").Value = -95120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 11933.871
$ws.Range("J4").Value = 11933.871
$ws.Range("L4").Value = 11933.871
$ws.Range("N4").Value = -12157.871
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = ""
$ws.Range("H52").Value = 79364
$ws.Range("J52").Value = 80348.39999999999
$ws.Range("L52").Value = 80348.39999999999
$ws.Range("N52").Value = -80936.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 37479536
$ws.Range("I4").Value = 59465330
$ws.Range("J4").Value = 103689.7
$ws.Range("K4").Value = 178395990
$ws.Range("L4").Value = 311069.1
$ws.Range("M4").Value = -178395878
$ws.Range("N4").Value = -311293.1
$ws.Range("H12").Value = 424.3125
$ws.Range("J12").Value = 514.7692
$ws.Range("L12").Value = 1544.3076
$ws.Range("N12").Value = -1890.3076

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = ""
$ws.Range("H45").Value = 63662.5
$ws.Range("J45").Value = 70000
$ws.Range("L45").Value = 70000
$ws.Range("N45").Value = -71118
$ws.Range("H135").Value = 95000
$ws.Range("J135").Value = 95000
$ws.Range("L135").Value = 95000
$ws.Range("N135").Value = -105140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2099.25
$ws.Range("I16").Value = 2099.25
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2099.25
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1929.25
$ws.Range("N16").Value = ""
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").Value = ""
$ws.Range("H46").Value = 5001423.5
$ws.Range("I46").Value = 10000885
$ws.Range("J46").Value = 1961.9
$ws.Range("K46").Value = 10000885
$ws.Range("L46").Value = 1961.9
$ws.Range("M46").Value = -10000697
$ws.Range("N46").Value = -2337.9
$ws.Range("H68").Value = 2121.5
$ws.Range("I68").Value = 2121.5
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2121.5
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1372.5
$ws.Range("N68").Value = ""
$ws.Range("H71").Value = 2121.5
$ws.Range("I71").Value = 2121.5
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 10607.5
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -6863.5
$ws.Range("N71").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 6999
$ws.Range("I2").Value = 6999
$ws.Range("K2").Value = 6999
$ws.Range("M2").Value = -6887
$ws.Range("H81").Value = 32680.625
$ws.Range("I81").Value = 11844.2
$ws.Range("J81").Value = 67408
$ws.Range("K81").Value = 23688.4
$ws.Range("L81").Value = 134816
$ws.Range("M81").Value = -22627.4
$ws.Range("N81").Value = -136938
$ws.Range("H84").Value = 32680.625
$ws.Range("I84").Value = 11844.2
$ws.Range("J84").Value = 67408
$ws.Range("K84").Value = 118442
$ws.Range("L84").Value = 674080
$ws.Range("M84").Value = -113138
$ws.Range("N84").Value = -684688
